# Insert a new data row at row 23 (pushes existing rows 23:117 down to 24:118,
# dimension grows from A1:T117 to A1:T118), then populate the new row 23 with
# the latest market observation for Granada (Vega Modelo de Temuco).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 23..117 down by one row.
$ws.Rows("23:23").Insert()

# New row 23 data.
$ws.Range("A23").Value = 10
$ws.Range("B23").Value = "Vega Modelo de Temuco"
$ws.Range("C23").Value = "La Araucanía"
$ws.Range("D23").Value = 44687
$ws.Range("E23").Value = 9
$ws.Range("F23").Value = "Fruta"
$ws.Range("G23").Value = 100104
$ws.Range("H23").Value = "Frutos de pepita"
$ws.Range("I23").Value = 100104001
$ws.Range("J23").Value = "Granada"
$ws.Range("K23").Value = "Wonderfull"
$ws.Range("L23").Value = "Primera"
$ws.Range("M23").Value = 65
$ws.Range("N23").Value = 14000
$ws.Range("O23").Value = 15000
$ws.Range("P23").Value = 14615
$ws.Range("Q23").Value = "`$/bandeja 15 kilos granel"
$ws.Range("R23").Value = "Provincia de Limarí"
$ws.Range("S23").Value = 974
$ws.Range("T23").Value = 15
